$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the boolean cell D1 entirely
$ws.Range("D1").Clear()

# Set raw numeric (date serial) values for A1, B1, C1
$ws.Range("A1").Value = 42894.35944359954
$ws.Range("B1").Value = 42894.35944395833
$ws.Range("C1").Value = 42894.35944395833

# Apply a date/time number format (built-in format 22: m/d/yy h:mm) to B1 and C1
$ws.Range("B1:C1").NumberFormat = "m/d/yy h:mm"
